$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "so_origin"
$ws.Range("A2").Value = "SO2385027"

$ws.Range("B7").Select()
